# Apply odds/score updates for Jogos_da_Semana_FlashScore_2025-02-12.xlsx
# Updates 206 numeric cells across rows 2-26 (Sheet1) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.5
$ws.Range("I2").Value = 1.42
$ws.Range("L2").Value = 1.91
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 67
$ws.Range("AM2").Value = 7.5

# Row 3
$ws.Range("G3").Value = 3.25
$ws.Range("I3").Value = 2.45
$ws.Range("J3").Value = 4.33
$ws.Range("L3").Value = 3.4
$ws.Range("AD3").Value = 41
$ws.Range("AO3").Value = 23
$ws.Range("AP3").Value = 26

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("AL5").Value = 7.5
$ws.Range("AN5").Value = 15

# Row 6
$ws.Range("G6").Value = 2.1
$ws.Range("I6").Value = 3.8
$ws.Range("J6").Value = 3
$ws.Range("W6").Value = 1.62
$ws.Range("X6").Value = 2.2
$ws.Range("Y6").Value = 2.25
$ws.Range("Z6").Value = 1.57
$ws.Range("AB6").Value = 8.5
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 19
$ws.Range("AG6").Value = 6
$ws.Range("AL6").Value = 8
$ws.Range("AM6").Value = 17

# Row 9
$ws.Range("L9").Value = 9.5
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 3.15
$ws.Range("T9").Value = 1.35
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 1.22
$ws.Range("AJ9").Value = 126
$ws.Range("AR9").Value = 1.61
$ws.Range("AS9").Value = 2.28

# Row 11
$ws.Range("G11").Value = 2.95
$ws.Range("H11").Value = 2.57
$ws.Range("I11").Value = 2.85
$ws.Range("J11").Value = 3.75
$ws.Range("K11").Value = 1.78
$ws.Range("L11").Value = 3.55
$ws.Range("M11").Value = 1.16
$ws.Range("N11").Value = 4.6
$ws.Range("W11").Value = 1.65
$ws.Range("X11").Value = 2.12
$ws.Range("Y11").Value = 2.15
$ws.Range("Z11").Value = 1.62
$ws.Range("AB11").Value = 13
$ws.Range("AD11").Value = 40
$ws.Range("AG11").Value = 4.6
$ws.Range("AH11").Value = 5.3
$ws.Range("AI11").Value = 18
$ws.Range("AJ11").Value = 120
$ws.Range("AL11").Value = 6.3
$ws.Range("AM11").Value = 13
$ws.Range("AO11").Value = 37
$ws.Range("AQ11").Value = 50

# Row 12
$ws.Range("G12").Value = 3.25
$ws.Range("H12").Value = 2.82
$ws.Range("J12").Value = 3.8
$ws.Range("L12").Value = 3.05
$ws.Range("N12").Value = 5.7
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.57
$ws.Range("Q12").Value = 2.3
$ws.Range("R12").Value = 1.55
$ws.Range("U12").Value = 4
$ws.Range("Y12").Value = 1.9
$ws.Range("AA12").Value = 8.25
$ws.Range("AC12").Value = 11.25
$ws.Range("AD12").Value = 45
$ws.Range("AE12").Value = 32
$ws.Range("AF12").Value = 40
$ws.Range("AG12").Value = 5.7
$ws.Range("AH12").Value = 5.5
$ws.Range("AM12").Value = 10.75
$ws.Range("AP12").Value = 23
$ws.Range("AQ12").Value = 37

# Row 13
$ws.Range("G13").Value = 2.18
$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 3.95
$ws.Range("J13").Value = 2.87
$ws.Range("K13").Value = 1.82
$ws.Range("L13").Value = 4.7
$ws.Range("M13").Value = 1.15
$ws.Range("N13").Value = 4.75
$ws.Range("O13").Value = 1.62
$ws.Range("P13").Value = 2.15
$ws.Range("Q13").Value = 2.8
$ws.Range("R13").Value = 1.38
$ws.Range("U13").Value = 5.1
$ws.Range("V13").Value = 1.13
$ws.Range("W13").Value = 1.62
$ws.Range("X13").Value = 2.18
$ws.Range("Y13").Value = 2.25
$ws.Range("Z13").Value = 1.57
$ws.Range("AA13").Value = 5.2
$ws.Range("AB13").Value = 8.75
$ws.Range("AD13").Value = 21
$ws.Range("AG13").Value = 4.75
$ws.Range("AH13").Value = 5.6
$ws.Range("AI13").Value = 20
$ws.Range("AL13").Value = 7.5
$ws.Range("AM13").Value = 19.5
$ws.Range("AN13").Value = 15
$ws.Range("AP13").Value = 55
$ws.Range("AQ13").Value = 75

# Row 14
$ws.Range("G14").Value = 1.27
$ws.Range("H14").Value = 4.85
$ws.Range("I14").Value = 11.5
$ws.Range("J14").Value = 1.75
$ws.Range("K14").Value = 2.3
$ws.Range("N14").Value = 7.5
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.3
$ws.Range("Q14").Value = 1.85
$ws.Range("R14").Value = 1.85
$ws.Range("W14").Value = 1.4
$ws.Range("X14").Value = 2.72
$ws.Range("Y14").Value = 2.5
$ws.Range("Z14").Value = 1.47
$ws.Range("AB14").Value = 5.1
$ws.Range("AC14").Value = 9.5
$ws.Range("AD14").Value = 6.9
$ws.Range("AF14").Value = 45
$ws.Range("AG14").Value = 7.5
$ws.Range("AH14").Value = 10
$ws.Range("AL14").Value = 23
$ws.Range("AM14").Value = 90

# Row 16
$ws.Range("H16").Value = 3.6
$ws.Range("J16").Value = 2.3
$ws.Range("K16").Value = 2.1
$ws.Range("U16").Value = 4
$ws.Range("V16").Value = 1.25
$ws.Range("Y16").Value = 2.1
$ws.Range("Z16").Value = 1.67
$ws.Range("AB16").Value = 7
$ws.Range("AC16").Value = 8.5
$ws.Range("AE16").Value = 15
$ws.Range("AG16").Value = 8
$ws.Range("AI16").Value = 19
$ws.Range("AJ16").Value = 67
$ws.Range("AK16").Value = 501
$ws.Range("AL16").Value = 12
$ws.Range("AN16").Value = 17
$ws.Range("AO16").Value = 51
$ws.Range("AP16").Value = 41

# Row 17
$ws.Range("M17").Value = 1.1
$ws.Range("N17").Value = 7

# Row 19
$ws.Range("G19").Value = 2.4
$ws.Range("I19").Value = 2.8
$ws.Range("J19").Value = 3.2
$ws.Range("Y19").Value = 1.83
$ws.Range("Z19").Value = 1.83
$ws.Range("AB19").Value = 11
$ws.Range("AC19").Value = 9.5
$ws.Range("AK19").Value = 301
$ws.Range("AN19").Value = 11

# Row 20
$ws.Range("Q20").Value = 1.95
$ws.Range("R20").Value = 1.9

# Row 22
$ws.Range("I22").Value = 2.8
$ws.Range("J22").Value = 3
$ws.Range("W22").Value = 1.36
$ws.Range("X22").Value = 3
$ws.Range("AD22").Value = 23

# Row 23
$ws.Range("K23").Value = 1.95
$ws.Range("L23").Value = 5
$ws.Range("M23").Value = 1.1
$ws.Range("N23").Value = 7
$ws.Range("R23").Value = 1.53
$ws.Range("Y23").Value = 2.2
$ws.Range("Z23").Value = 1.62
$ws.Range("AA23").Value = 5.5
$ws.Range("AF23").Value = 41
$ws.Range("AI23").Value = 21
$ws.Range("AJ23").Value = 81

# Row 24
$ws.Range("G24").Value = 1.67
$ws.Range("H24").Value = 3.7
$ws.Range("I24").Value = 5
$ws.Range("J24").Value = 2.3
$ws.Range("K24").Value = 2.2
$ws.Range("L24").Value = 5
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10
$ws.Range("Q24").Value = 1.93
$ws.Range("R24").Value = 1.88
$ws.Range("W24").Value = 1.4
$ws.Range("X24").Value = 2.75
$ws.Range("AB24").Value = 8
$ws.Range("AD24").Value = 13
$ws.Range("AG24").Value = 10
$ws.Range("AH24").Value = 7

# Row 25
$ws.Range("M25").Value = 1.07
$ws.Range("N25").Value = 9
$ws.Range("AQ25").Value = 34

# Row 26
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
